$p = $ppt.ActivePresentation
$s = $p.Slides.Add(4, 7)
$sh = $s.Shapes.AddShape(9, 100, 100, 200, 50)
$sh.Name = "円/楕円 16"
Write-Host "ok"
